$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 2: B2 77719 -> 77733
$ws.Range("B2").Value = 77733

# Row 3: B3 73820 -> 73834
$ws.Range("B3").Value = 73834

# Row 4: B4 85836 -> 85850
$ws.Range("B4").Value = 85850

# Row 5: swap taxon info with row 6's old values, plus new Id/B values
$ws.Range("A5").Value = 112231361
$ws.Range("B5").Value = 77650
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."

# Row 6: swap taxon info with row 5's old values, plus new Id/B values
$ws.Range("A6").Value = 112231346
$ws.Range("B6").Value = 76634
$ws.Range("E6").Value = 228579
$ws.Range("F6").Value = "Liten svartspik"
$ws.Range("G6").Value = "Chaenothecopsis nana"
$ws.Range("H6").Value = "Tibell"
